$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44216) -> becomes the data previously on row 6 (date 44253)
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 90
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12667
$ws.Range("S2").Value = 905

# Row 3 (was date 44232) -> becomes the data previously on row 7 (date 44229)
$ws.Range("D3").Value = 44229
$ws.Range("M3").Value = 55
$ws.Range("P3").Value = 11364
$ws.Range("S3").Value = 812

# Row 6 (was date 44253) -> becomes the data previously on row 2 (date 44216)
$ws.Range("D6").Value = 44216
$ws.Range("M6").Value = 55
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11545
$ws.Range("S6").Value = 825

# Row 7 (was date 44229) -> becomes the data previously on row 8 (date 44181)
$ws.Range("D7").Value = 44181
$ws.Range("M7").Value = 65
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9462
$ws.Range("S7").Value = 676

# Row 8 (was date 44181) -> becomes the data previously on row 3 (date 44232)
$ws.Range("D8").Value = 44232
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11583
$ws.Range("S8").Value = 827
